$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'basketball under pants'
$ws.Cells.Item(2, 1).Value = 'softball gear for girls'
$ws.Cells.Item(3, 1).Value = 'running capri'
$ws.Cells.Item(4, 1).Value = 'softball compression sleeve'
$ws.Cells.Item(5, 1).Value = 'youth softball compression sleeve'
$ws.Cells.Item(6, 1).Value = 'running tights mens'
$ws.Cells.Item(7, 1).Value = 'spandex men'
$ws.Cells.Item(8, 1).Value = 'hockey kneepads'
$ws.Cells.Item(9, 1).Value = 'padded leg sleeve'
$ws.Cells.Item(10, 1).Value = 'mens basketball gear'
$ws.Cells.Item(11, 1).Value = 'snowboarding padded shorts'
$ws.Cells.Item(12, 1).Value = 'padded shorts snowboarding'
$ws.Cells.Item(13, 1).Value = 'knee sleeve wrestling'
$ws.Cells.Item(14, 1).Value = 'sleeve knee pads'
$ws.Cells.Item(15, 1).Value = 'womens compression leggings'
$ws.Cells.Item(16, 1).Value = 'airsoft knee pads'
$ws.Cells.Item(17, 1).Value = 'mens compression tights 3 4'
$ws.Cells.Item(18, 1).Value = 'basketball clothes for men'
$ws.Cells.Item(19, 1).Value = 'men running tights'
$ws.Cells.Item(20, 1).Value = 'knee pads nike'
$ws.Cells.Item(21, 1).Value = 'knee pads mizuno'
$ws.Cells.Item(22, 1).Value = 'knee pads bike'
$ws.Cells.Item(23, 1).Value = 'yoga capri pants'
$ws.Cells.Item(24, 1).Value = 'knee pads mtb'
$ws.Cells.Item(25, 1).Value = 'knee pads skating'
$ws.Cells.Item(26, 1).Value = 'mens workout tights'
$ws.Cells.Item(27, 1).Value = 'mens basketball pants'
$ws.Cells.Item(28, 1).Value = 'asics knee pads'
$ws.Cells.Item(29, 1).Value = 'mens workout tights pants'
$ws.Cells.Item(30, 1).Value = 'downhill knee pads'
$ws.Cells.Item(31, 1).Value = 'men gym pants'
$ws.Cells.Item(32, 1).Value = 'athletic capris'
$ws.Cells.Item(33, 1).Value = 'valken knee pads'
$ws.Cells.Item(34, 1).Value = 'woodland knee pads'
$ws.Cells.Item(35, 1).Value = 'training tights men'
$ws.Cells.Item(36, 1).Value = 'short tights for men'
$ws.Cells.Item(37, 1).Value = 'ua compression pants'
$ws.Cells.Item(38, 1).Value = 'men workout tights'
$ws.Cells.Item(39, 1).Value = 'knee pads for exercise'
$ws.Cells.Item(40, 1).Value = 'mens leggins'
$ws.Cells.Item(41, 1).Value = 'nike kneepads'
$ws.Cells.Item(42, 1).Value = 'youth football girdle with knee pads'
$ws.Cells.Item(43, 1).Value = 'compression tights with pads'
$ws.Cells.Item(44, 1).Value = 'compression leggings with knee pads'
$ws.Cells.Item(45, 1).Value = 'basketball padded compression pants'
$ws.Cells.Item(46, 1).Value = 'basketball knee pad pants'
$ws.Cells.Item(47, 1).Value = 'padded compression pants men basketball'
$ws.Cells.Item(48, 1).Value = 'padded tights men basketball'
$ws.Cells.Item(49, 1).Value = 'tights with pads basketball'
$ws.Cells.Item(50, 1).Value = 'mens leggings with knee pads'
$ws.Cells.Item(51, 1).Value = 'leggings with knee pads women'
$ws.Cells.Item(52, 1).Value = 'mtb knee pads men'
$ws.Cells.Item(53, 1).Value = 'basketball padded knee sleeve'
$ws.Cells.Item(54, 1).Value = 'compression pants women'
$ws.Cells.Item(55, 1).Value = 'compression knee sleeve men basketball'
$ws.Cells.Item(56, 1).Value = 'basketball sweat pants for men'
$ws.Cells.Item(57, 1).Value = 'knee sleeve for wrestling'
$ws.Cells.Item(58, 1).Value = 'leg sleeves for basketball youth'
$ws.Cells.Item(59, 1).Value = 'training pants men'
$ws.Cells.Item(60, 1).Value = 'compression knee sleeve men basketball'
$ws.Cells.Item(61, 1).Value = 'basketball sweat pants for men'
$ws.Cells.Item(62, 1).Value = 'knee sleeve for wrestling'
$ws.Cells.Item(63, 1).Value = 'leg sleeves for basketball youth'
$ws.Cells.Item(64, 1).Value = 'goalkeeper knee pads'
$ws.Cells.Item(65, 1).Value = 'basketball calf sleeve'
$ws.Cells.Item(66, 1).Value = 'compression knee sleeves with pads'
$ws.Cells.Item(67, 1).Value = 'compression sleeve knee pads'
$ws.Cells.Item(68, 1).Value = 'youth knee sleeve'
$ws.Cells.Item(69, 1).Value = 'knee pad for scooter'
$ws.Cells.Item(70, 1).Value = 'knee basketball'
$ws.Cells.Item(71, 1).Value = 'knee pads for basketball youth'
$ws.Cells.Item(72, 1).Value = 'mens compression knee'
$ws.Cells.Item(73, 1).Value = 'knee pad sleeve basketball'
$ws.Cells.Item(74, 1).Value = 'mens 3/4 compression pants'
$ws.Cells.Item(75, 1).Value = 'youth compression knee pad sleeve'
$ws.Cells.Item(76, 1).Value = 'mens basketball knee sleeves'
$ws.Cells.Item(77, 1).Value = 'knee sleeve wrestling youth'
$ws.Cells.Item(78, 1).Value = 'knee sleeves basketball men'
$ws.Cells.Item(79, 1).Value = 'soccer compression pants'
$ws.Cells.Item(80, 1).Value = 'leggings tight'
$ws.Cells.Item(81, 1).Value = 'basketball leg sleeve youth padded'
$ws.Cells.Item(82, 1).Value = 'knee pad construction'
$ws.Cells.Item(83, 1).Value = 'youth basketball knee sleeve'
$ws.Cells.Item(84, 1).Value = 'working knee pads for men'
$ws.Cells.Item(85, 1).Value = 'cycling pants for men'
$ws.Cells.Item(86, 1).Value = 'boys youth leggings'
$ws.Cells.Item(87, 1).Value = 'compression running capris'
$ws.Cells.Item(88, 1).Value = 'knee sleeve baseball'
$ws.Cells.Item(89, 1).Value = 'compression knee sleeves for basketball'
$ws.Cells.Item(90, 1).Value = 'volleyball kneepads'
$ws.Cells.Item(91, 1).Value = 'compression knee sleeve with pad'
$ws.Cells.Item(92, 1).Value = 'men capri shorts'
$ws.Cells.Item(93, 1).Value = 'running compression pants'
$ws.Cells.Item(94, 1).Value = 'mens work pants knee pads'
$ws.Cells.Item(95, 1).Value = 'best knee pads'
$ws.Cells.Item(96, 1).Value = 'compression pants sleeves'
$ws.Cells.Item(97, 1).Value = 'mens compression running tights'
$ws.Cells.Item(98, 1).Value = 'knee pads working'
$ws.Cells.Item(99, 1).Value = 'basketball aids'
$ws.Cells.Item(100, 1).Value = 'baseball youth pants'
